$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()

$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value = "Los Lagos"
$ws.Cells.Item(17, 4).Value = 44533
$ws.Cells.Item(17, 5).Value = 10
$ws.Cells.Item(17, 6).Value = 300000000
$ws.Cells.Item(17, 7).Value = "Espárragos"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 420
$ws.Cells.Item(17, 11).Value = 1700
$ws.Cells.Item(17, 12).Value = 1700
$ws.Cells.Item(17, 13).Value = 1700
$ws.Cells.Item(17, 14).Value = "$/kilo"
$ws.Cells.Item(17, 15).Value = "Provincia de Linares"
$ws.Cells.Item(17, 16).Value = 1700
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"
